$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Waktu Pencarian")
$ws.Select()

$rng = $ws.Range("A1:G65")
$sortField = $ws.Range("G1:G65")

$rng.Sort($sortField, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 2, [System.Reflection.Missing]::Value, 1, 1)

$ws.Range("K8").Select()
